$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @(Price, Volume(1h)) ; Price values that parse as plain numbers are
# apostrophe-prefixed so Excel stores them as text (matching the sheet's original
# inline-string cells) instead of silently converting them to numeric values.
$updates = @{
    2 = @('30.452.07', '  +0.34%  ')
    3 = @('1.941.90', '  +0.31%  ')
    4 = @('''1.006', '  +0.55%  ')
    5 = @('''0.7551', '  +4.02%  ')
    6 = @('''246.71', '  -2.17%  ')
    7 = @('''1.006', '  +0.59%  ')
    8 = @('''0.3194', '  -3.41%  ')
    9 = @('''27.69', '  -0.97%  ')
    10 = @('''0.07009', '  -3.31%  ')
    11 = @('''0.7820', '  -3.37%  ')
    12 = @('''0.08023', '  -1.02%  ')
    13 = @('1.934.89', '  -0.01%  ')
    14 = @('''5.354', '  -2.26%  ')
    15 = @('''94.58', '  -0.30%  ')
    16 = @('''14.45', '  -4.37%  ')
    17 = @('30.430.54', '  +0.28%  ')
    18 = @('''254.58', '  +0.55%  ')
    19 = @('''0.000007928', '  -3.76%  ')
    20 = @('''5.776', '  -0.88%  ')
    21 = @('2.187.13', '  -0.15%  ')
    22 = @('''1.004', '  +0.41%  ')
    23 = @('''1.007', '  +0.67%  ')
    24 = @('''6.678', '  -4.10%  ')
    25 = @('''9.527', '  -2.42%  ')
    26 = @('''165.31', '  -0.47%  ')
    27 = @('''19.05', '  -1.49%  ')
    28 = @('''0.1330', '  +2.27%  ')
    29 = @('''2.268', '  -3.44%  ')
    30 = @('''1.377', '  +1.73%  ')
    31 = @('''1.516', '  -2.03%  ')
    32 = @('''4.394', '  -1.05%  ')
    33 = @('''4.119', '  -2.29%  ')
    34 = @('''0.05158', '  -1.68%  ')
    35 = @('''1.279', '  +0.83%  ')
    36 = @('''0.7458', '  -0.68%  ')
    37 = @('''2.797', '  +0.74%  ')
    38 = @('''0.01947', '  -1.25%  ')
    39 = @('''2.814', '  +0.40%  ')
    40 = @('''78.88', '  -0.70%  ')
    41 = @('''6.415', '  -0.58%  ')
    42 = @('''0.4487', '  -1.42%  ')
    43 = @('''1.966', '  -3.26%  ')
    44 = @('''1.006', '  +0.53%  ')
    45 = @('''0.8330', '  -1.40%  ')
    46 = @('''101.19', '  -0.85%  ')
    47 = @('''9.768', '  -0.58%  ')
    48 = @('''7.481', '  +0.40%  ')
    49 = @('''37.20', '  +1.09%  ')
    50 = @('''976.40', '  +9.93%  ')
    51 = @('''0.06030', '  -0.33%  ')
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]   # column D = Price
    $ws.Cells.Item($row, 5).Value = $vals[1]   # column E = Volume(1h)
}
